$wb = $excel.ActiveWorkbook

# ===== Sheet: 展览 (numeric 'F' column updates) =====
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 57
$ws1.Range("F3").Value = 7455
$ws1.Range("F4").Value = 3551
$ws1.Range("F6").Value = 3875
$ws1.Range("F7").Value = 70
$ws1.Range("F8").Value = 92
$ws1.Range("F10").Value = 112
$ws1.Range("F11").Value = 164
$ws1.Range("F12").Value = 517
$ws1.Range("F13").Value = 19
$ws1.Range("F14").Value = 157
$ws1.Range("F17").Value = 354
$ws1.Range("F18").Value = 4195
$ws1.Range("F21").Value = 1034
$ws1.Range("F22").Value = 537
$ws1.Range("F23").Value = 1903
$ws1.Range("F25").Value = 102
$ws1.Range("F26").Value = 58
$ws1.Range("F27").Value = 3079
$ws1.Range("F28").Value = 2321
$ws1.Range("F29").Value = 69
$ws1.Range("F30").Value = 86
$ws1.Range("F32").Value = 82
$ws1.Range("F33").Value = 124
$ws1.Range("F36").Value = 110
$ws1.Range("F37").Value = 4411
$ws1.Range("F38").Value = 504
$ws1.Range("F39").Value = 327
$ws1.Range("F42").Value = 836
$ws1.Range("F43").Value = 236
$ws1.Range("F45").Value = 1665
$ws1.Range("F46").Value = 262

# ----- Sheet 展览: Row 35 full content replace (event swapped) -----
$ws1.Range("C35").Value = '北京·广播剧《踏雪乌啼观海啸》·专场活动'
$ws1.Range("D35").Value = '京沈路与天北路交汇处西北角 中国国际展览中心新馆'
$ws1.Range("E35").Value = '2024.07.21 12:00-07.21 15:15'
$ws1.Range("F35").Value = 2
$ws1.Range("G35").Value = 288
$ws1.Range("H35").Value = 'https://show.bilibili.com/platform/detail.html?id=86917'
$ws1.Range("I35").Value = '//i1.hdslb.com/bfs/openplatform/202406/0yWKaRsR1717580159808.png'

# ===== Sheet: 演出 (numeric updates + new row insert) =====
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 3
$ws2.Range("F4").Value = 443
$ws2.Range("F7").Value = 63
$ws2.Range("F13").Value = 1

# Insert a new row at position 15 (shifts old rows 15-20 down to 16-21)
$ws2.Rows.Item(15).Insert()
# Copy formatting (style) of A column from the row below into the new A15 cell
$ws2.Range("A16").Copy()
$ws2.Range("A15").PasteSpecial(-4122)
$ws2.Range("A15").Value = 14

# New row 15 content (brand-new event)
$ws2.Range("B15").Value = '2024-07-23'
$ws2.Range("C15").Value = '北京·巴西浪漫风情——手风琴大满贯音乐家道格拉斯·博尔萨蒂专场音乐会'
$ws2.Range("D15").Value = '复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)'
$ws2.Range("E15").Value = '2024.07.23 19:30-07.23 21:00'
$ws2.Range("F15").Value = 0
$ws2.Range("G15").Value = 140
$ws2.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=86922'
$ws2.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202405/i14RABlz1716527544509.jpeg'

# Fix F19 (old row18's 'Marcin Patrzalek' event, now shifted to row 19): 604 -> 608
$ws2.Range("F19").Value = 608

# ===== Sheet: 全部类型 (numeric 'F' column updates) =====
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3
$ws4.Range("F5").Value = 7455
$ws4.Range("F6").Value = 3551
$ws4.Range("F7").Value = 3875
$ws4.Range("F8").Value = 70
$ws4.Range("F9").Value = 92
$ws4.Range("F11").Value = 112
$ws4.Range("F13").Value = 164
$ws4.Range("F14").Value = 517
$ws4.Range("F16").Value = 157
$ws4.Range("F19").Value = 354
$ws4.Range("F20").Value = 4195
$ws4.Range("F26").Value = 537
$ws4.Range("F27").Value = 1903
$ws4.Range("F29").Value = 102
$ws4.Range("F30").Value = 3079
$ws4.Range("F31").Value = 2322
$ws4.Range("F32").Value = 69
$ws4.Range("F33").Value = 86
$ws4.Range("F34").Value = 124
$ws4.Range("F35").Value = 110
$ws4.Range("F37").Value = 4411
$ws4.Range("F39").Value = 504
$ws4.Range("F40").Value = 327
$ws4.Range("F42").Value = 836
$ws4.Range("F43").Value = 236
$ws4.Range("F45").Value = 1665
$ws4.Range("F46").Value = 262

$wb.Save()
